$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44322
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("S2").Value = 1714

$ws.Range("D3").Value = 44322
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("S3").Value = 1143

$ws.Range("D4").Value = 44300
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 2143

$ws.Range("D5").Value = 44300
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("S5").Value = 1714

$ws.Range("D8").Value = 44301
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("S8").Value = 2000

$ws.Range("D9").Value = 44301
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

$ws.Range("D10").Value = 44292
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("S10").Value = 2286

$ws.Range("D11").Value = 44292
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 2143

$ws.Range("D12").Value = 44302
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 2143

$ws.Range("D13").Value = 44302
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 1714
